$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: add Date/Time, Topics, Participants (ID already present in A6, Participants-summary already in E6)
$ws.Range("B2").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("B6").Value = 45659
$ws.Range("C6").Value = "Discussed current progress on Chat Profiles, divided up the work, and planned for the report writing."
$ws.Range("D6").Value = "ma, is, se, cl, ce"

# Row 7: add Date/Time, Topics, Participants (ID already present in A7)
$ws.Range("B2").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B7").Value = 45664
$ws.Range("C7").Value = "Reviewed current work, Suggest Updates to One of the Chats, Planned Remaining Report Writing"
$ws.Range("D7").Value = "ma, is, se, cl, ce"

$excel.CutCopyMode = 0

# Update selection to C6 (matches the saved cursor position in the diff)
$ws.Range("C6").Select()
